$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 117, shifting existing rows 117:221 down to 118:222
$ws.Rows.Item(117).Insert()

# Populate the newly inserted row 117 with the new record's data
$ws.Cells.Item(117, 1).Value = 1
$ws.Cells.Item(117, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(117, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(117, 4).Value = 44669
$ws.Cells.Item(117, 5).Value = 15
$ws.Cells.Item(117, 6).Value = "Fruta"
$ws.Cells.Item(117, 7).Value = 100108
$ws.Cells.Item(117, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(117, 9).Value = 100108006
$ws.Cells.Item(117, 10).Value = "Plátano"
$ws.Cells.Item(117, 11).Value = "Sin especificar"
$ws.Cells.Item(117, 12).Value = "Pintón"
$ws.Cells.Item(117, 13).Value = 120
$ws.Cells.Item(117, 14).Value = 21000
$ws.Cells.Item(117, 15).Value = 22000
$ws.Cells.Item(117, 16).Value = 21500
$ws.Cells.Item(117, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(117, 18).Value = "Ecuador"
$ws.Cells.Item(117, 19).Value = 1075
$ws.Cells.Item(117, 20).Value = 20
